# Apply v0.1.2 update:
#   1. Bump "Version" value on the Metadata sheet from 0.1.1 -> 0.1.2
#   2. Add a new "Include from VSTiposDocumento" sheet (clone of the
#      existing "Include from Tipo Identificad" include-sheet) whose
#      System URI points at the VSTiposDocumentos value set.

$wb = $excel.ActiveWorkbook

# --- 1. Metadata!Version 0.1.1 -> 0.1.2 --------------------------------
$metadata = $wb.Worksheets.Item(1)
[void]$metadata.Cells.Replace("0.1.1", "0.1.2")

# --- 2. New "Include from VSTiposDocumento" sheet ----------------------
$template = $wb.Worksheets.Item(2)
$lastIndex = $wb.Worksheets.Count
$newSheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item($lastIndex))
$newSheet.Name = "Include from VSTiposDocumento"

# Clone the template sheet's layout/content/formatting cell by cell so
# the new include-sheet matches the existing one's styling.
$template.Range("A1").Copy($newSheet.Range("A1"))
$template.Range("A2").Copy($newSheet.Range("A2"))
$template.Range("A3:B3").Copy($newSheet.Range("A3"))
$template.Range("A4:B4").Copy($newSheet.Range("A4"))

# Column widths matching the template sheet.
$newSheet.Columns.Item(1).ColumnWidth = 30.703125
$newSheet.Columns.Item(2).ColumnWidth = 50.703125

# Point the new sheet's System URI at the VSTiposDocumentos value set.
$newSheet.Range("B4").Value = "https://hl7chile.cl/fhir/ig/CoreCL/ValueSet/VSTiposDocumentos"

# Restore original active sheet/selection.
$metadata.Activate()
